# Adapt the DE results table to the new PSO ("MaxFES") layout:
#  - Column A header "Gen" -> "MaxFES", and its generation-count values
#    replaced by the fraction-of-budget values used by the PSO run.
#  - Drop the old "Run 50" column (AZ) entirely.
#  - The old "Mean" column (BA) slides left into AZ and is recomputed
#    over only the 50 remaining run columns (B:AY).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: header + MaxFES values -------------------------------------
$ws.Range("A1").Value = "MaxFES"

$maxFes = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
$row = 2
foreach ($v in $maxFes) {
    $ws.Cells.Item($row, 1).Value = $v
    $row = $row + 1
}

# --- Recompute the Mean column (currently AZ = old "Run 50") over the ------
# --- remaining 50 run columns (B:AY) before the old Mean column (BA) is ----
# --- deleted and this one shifts into its place. ---------------------------
$meanValues = @{
    2  = 151119993710.4635
    3  = 150263291216.7626
    4  = 120110881240.5643
    5  = 38209799619.06186
    6  = 12276547640.18491
    7  = 3118448940.396138
    8  = 700114027.4248766
    9  = 159664737.0636565
    10 = 40111430.33000951
    11 = 9942609.65779773
    12 = 2384516.60491061
    13 = 681309.34633288
    14 = 172439.3936323
}

foreach ($r in $meanValues.Keys) {
    $ws.Cells.Item($r, 52).Value = $meanValues[$r]
}

$ws.Cells.Item(1, 52).Value = "Mean"

# --- Drop the old "Run 50" data (now superseded) by removing the old ------
# --- trailing "Mean" column BA; AZ above already holds the refreshed ------
# --- mean so the sheet ends up exactly one column narrower (A:AZ). --------
$ws.Columns("BA").Delete()
